$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "does user have coupon" input column (E) ---------------------------
# Fill column E top-to-bottom first so new shared-string entries come out in
# the same order as the authored workbook (yes/no before green/blue/purple).
$ws.Range("E8").Value = "input: does user have coupon"
$ws.Range("E9").Value = "yes"
$ws.Range("E10").Value = "no"
$ws.Range("E11").Value = "no"
$ws.Range("E12").Value = "yes"

# --- Package-type column (B) -------------------------------------------------
$ws.Range("B9").Value = "green"
$ws.Range("B10").Value = "blue"
$ws.Range("B11").Value = "purple"
$ws.Range("B12").Value = "green"

# --- Months purchasing for (C) / Additional data needed (D) ----------------
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 0

$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 0

$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 0

$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 5

# --- Output formulas ---------------------------------------------------------
# Row 9 (Green, 2 months, no extra data) - standalone formulas.
$ws.Range("F9").Formula = "=C2*C9"
$ws.Range("G9").Formula = "=C2*C9-20"

# Rows 10-11 (Blue / Purple, 2 months) share one fill-down formula, entered as
# a single range assignment so Excel records it as a shared formula group.
$ws.Range("F10:F11").Formula = "=C3*C10"
$ws.Range("G10:G11").Formula = "=C3*C10"

# Row 12 (Green, 3 months + 5 GB extra data) needs its own formula shape.
$ws.Range("F12").Formula = "=(C2*C12)+(D12*D2)"
$ws.Range("G12").Formula = "=(C2*C12)+(D12*D2)-20"

# --- Row heights (auto height shrank once the wrapped text was finalized) --
$ws.Rows.Item(9).RowHeight = 31
$ws.Rows.Item(10).RowHeight = 31
$ws.Rows.Item(11).RowHeight = 31

# --- Selection moved to G13 before the file was saved -----------------------
$ws.Range("G13").Select()
